$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The H2:H43 column currently holds "=TRUE()" boolean formulas (numeric
# result 1). The fix replaces them with the literal text value "TRUE"
# (a plain string, not a boolean) - matching the commit's "boolean values"
# fix where is_active became a text column instead of a computed boolean.
#
# A direct `.Value = "TRUE"` assignment gets auto-coerced back into an
# Excel boolean (t="b"), so instead we stage the literal text in a scratch
# range (using a leading apostrophe to force text interpretation), copy it,
# and paste-special just the values into H2:H43 - this carries over the
# text type without disturbing the existing cell style.
$stage = $ws.Range("Z1:Z42")
$stage.Value = "'TRUE"
$stage.Copy()
$ws.Range("H2:H43").PasteSpecial(-4163)  # xlPasteValues

# Remove the scratch column so it doesn't linger in the saved sheet / widen
# the used range.
$ws.Columns.Item(26).Delete()

# Match the updated view state: selection now spans the whole is_active
# column of data (H2:H43) instead of the header-only H1 selection.
$ws.Range("H2:H43").Select() | Out-Null
